$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attribute name, datatype, and display name for row 2
$ws.Range("B2").Value = "REVENUE_STATUS"
$ws.Range("C2").Value = "mdex:string"
$ws.Range("E2").Value = "Revenue Status"

# Update the active selection to C2 (as reflected in the sheet view)
$ws.Range("C2").Select()
